# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list (under the "Impact"
# sub-heading) so it reads as concise, impact-focused accomplishment
# statements instead of the longer job-duty style bullets, trimming the
# list from six items down to four.
#
# Original six bullets, in document order (under KEY ACHIEVEMENTS AND IMPACT
# / Impact):
#   1. Delivered $4.9M additional revenue through continuous testing and
#      optimization, increased conversion rates by 23%
#   2. Built redistricting platform used by thousands of analysts
#      nationwide with real-time collaborative editing and Census
#      integration, serving 12,847 analysts across 89 organizations
#   3. Achieved 87% prediction accuracy for voter turnout vs. industry
#      standard of 71%, reducing polling error margins from +/-4.2% to
#      +/-2.1%
#   4. Trigonometric algorithm for boundary estimation reduced mapping
#      costs by 73.5%, saving campaigns and organizations $4.7M and
#      enabling smaller nonprofits to conduct analysis        <- removed
#   5. Discovered systematic race coding errors affecting all Black and
#      Asian-American voters, developed geospatial machine learning
#      algorithms improving classification accuracy from 23% to 64%  <- removed
#   6. Developed longitudinal data analysis methods using geospatial
#      techniques that improved segmentation accuracy by 34% and survey
#      incidence rates by 28%, reducing polling costs while increasing
#      response quality
#
# New four bullets:
#   1. Platform impact: Built redistricting system serving 12,847
#      analysts across 89 organizations
#   2. Real-time collaboration at national scale
#   3. Revenue generation: Delivered $4.9M additional revenue through
#      optimization
#   4. 23% conversion rate improvement
#
# Several of the bullet sentences (e.g. "Achieved 87% prediction accuracy
# ...") also appear verbatim elsewhere in the resume (Professional
# Experience), so this script locates the "KEY ACHIEVEMENTS AND IMPACT"
# section by its heading and edits paragraphs by position within that
# section rather than relying on a document-wide text search, to avoid
# touching the look-alike bullet under Professional Experience.

$d = $word.ActiveDocument

# Find the paragraph index of the "KEY ACHIEVEMENTS AND IMPACT" heading.
$headingIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.TrimEnd("`r", "`n")
    if ($t -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -gt 0) {
    # Layout relative to the heading:
    #   headingIndex + 0 : "KEY ACHIEVEMENTS AND IMPACT" (Heading2)
    #   headingIndex + 1 : "Impact" (Heading3)
    #   headingIndex + 2 : bullet 1 (Delivered $4.9M ...)
    #   headingIndex + 3 : bullet 2 (Built redistricting platform ...)
    #   headingIndex + 4 : bullet 3 (Achieved 87% ...)
    #   headingIndex + 5 : bullet 4 (Trigonometric algorithm ...)  -> delete
    #   headingIndex + 6 : bullet 5 (Discovered systematic race ...) -> delete
    #   headingIndex + 7 : bullet 6 (Developed longitudinal ...)

    # Delete bullets 4 and 5 (indices headingIndex+5 and headingIndex+6),
    # highest index first so the other index stays valid.
    $d.Paragraphs($headingIndex + 6).Range.Delete()
    $d.Paragraphs($headingIndex + 5).Range.Delete()

    # Rewrite the remaining four bullets (now contiguous) with the new,
    # shorter accomplishment-focused wording, preserving each paragraph's
    # own formatting by writing through Range.Text.
    $d.Paragraphs($headingIndex + 2).Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
    $d.Paragraphs($headingIndex + 3).Range.Text = "• Real-time collaboration at national scale"
    $d.Paragraphs($headingIndex + 4).Range.Text = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
    $d.Paragraphs($headingIndex + 5).Range.Text = "• 23% conversion rate improvement"
}

Write-Output "Key Achievements section rewritten."
